# Weekly price update: insert two new rows of fresh "Cebolla" price data at
# the top of the existing date-ordered block (rows 706-803), pushing the
# existing rows down by two (to 708-805).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 706-707; everything from old row 706 down
# shifts to 708 onward (Excel carries the row-706 formatting, incl. the
# date number-format on column D, onto the freshly inserted rows).
$ws.Range("A706:A707").EntireRow.Insert()

# New row 706: "1a (guarda)" lot
$ws.Cells.Item(706, 1).Value  = 8
$ws.Cells.Item(706, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(706, 3).Value  = "Coquimbo"
$ws.Cells.Item(706, 4).Value  = 44776
$ws.Cells.Item(706, 5).Value  = 4
$ws.Cells.Item(706, 6).Value  = 100112004
$ws.Cells.Item(706, 7).Value  = "Cebolla"
$ws.Cells.Item(706, 8).Value  = "Sin especificar"
$ws.Cells.Item(706, 9).Value  = "1a (guarda)"
$ws.Cells.Item(706, 10).Value = 2800
$ws.Cells.Item(706, 11).Value = 7400
$ws.Cells.Item(706, 12).Value = 7500
$ws.Cells.Item(706, 13).Value = 7450
$ws.Cells.Item(706, 14).Value = "$/malla 16 kilos"
$ws.Cells.Item(706, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(706, 16).Value = 466
$ws.Cells.Item(706, 17).Value = 16
$ws.Cells.Item(706, 18).Value = "Hortaliza"

# New row 707: "2a (guarda)" lot
$ws.Cells.Item(707, 1).Value  = 8
$ws.Cells.Item(707, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(707, 3).Value  = "Coquimbo"
$ws.Cells.Item(707, 4).Value  = 44776
$ws.Cells.Item(707, 5).Value  = 4
$ws.Cells.Item(707, 6).Value  = 100112004
$ws.Cells.Item(707, 7).Value  = "Cebolla"
$ws.Cells.Item(707, 8).Value  = "Sin especificar"
$ws.Cells.Item(707, 9).Value  = "2a (guarda)"
$ws.Cells.Item(707, 10).Value = 1760
$ws.Cells.Item(707, 11).Value = 7000
$ws.Cells.Item(707, 12).Value = 7200
$ws.Cells.Item(707, 13).Value = 7100
$ws.Cells.Item(707, 14).Value = "$/malla 16 kilos"
$ws.Cells.Item(707, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(707, 16).Value = 444
$ws.Cells.Item(707, 17).Value = 16
$ws.Cells.Item(707, 18).Value = "Hortaliza"
